# Weekly Fruit/Vegetable price update:
# A new observation is inserted at row 70 (pushing all existing rows 70-173
# down by one, growing the table to row 174), and the new row 70 is
# populated with a fresh "Acelga" price record.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 70..173 down to 71..174, carrying along all formatting
# (including the date-number style on column D).
$ws.Rows.Item(70).Insert()

# Populate the newly-inserted row 70 with the new record.
$ws.Cells.Item(70, 1).Value2 = 7
$ws.Cells.Item(70, 2).Value2 = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(70, 3).Value2 = 'Ñuble'
$ws.Cells.Item(70, 4).Value2 = 44495
$ws.Cells.Item(70, 5).Value2 = 16
$ws.Cells.Item(70, 6).Value2 = 100112009
$ws.Cells.Item(70, 7).Value2 = 'Acelga'
$ws.Cells.Item(70, 8).Value2 = 'Sin especificar'
$ws.Cells.Item(70, 9).Value2 = 'Primera'
$ws.Cells.Item(70, 10).Value2 = 60
$ws.Cells.Item(70, 11).Value2 = 350
$ws.Cells.Item(70, 12).Value2 = 400
$ws.Cells.Item(70, 13).Value2 = 375
$ws.Cells.Item(70, 14).Value2 = '$/atado 0,5 a 1 kilo'
$ws.Cells.Item(70, 15).Value2 = 'Provincia de Diguillín'
$ws.Cells.Item(70, 16).Value2 = 375
$ws.Cells.Item(70, 17).Value2 = 1
$ws.Cells.Item(70, 18).Value2 = 'Hortaliza'
